# Automatische test-sync: 2025-08-03 23:26:50
# Append a new log row (row 50) to the "Logs" sheet, extend the
# conditional-formatting ranges that tracked the old last row (49) to
# the new last row (50), and bump the "Planning / Afspraak" tally on
# the "Dashboard" sheet from 13 to 14.

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")

$newRow = 50
$logs.Cells.Item($newRow, 1).Value = "Kun jij dit afhandelen?"
$logs.Cells.Item($newRow, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item($newRow, 3).Value = "Testmail #1: Kun jij dit afhandelen?"
$logs.Cells.Item($newRow, 4).Value = "Planning / Afspraak"
$logs.Cells.Item($newRow, 5).Value = "Bedankt, we hebben dit doorgestuurd naar planning@bedrijf.nl."
$logs.Cells.Item($newRow, 6).Value = "2025-08-03 23:26:42"
$logs.Cells.Item($newRow, 7).Value = "Ja"
$logs.Cells.Item($newRow, 8).Value = "Ja"
$logs.Cells.Item($newRow, 9).Value = "Nee"
$logs.Cells.Item($newRow, 10).Value = "Nee"

# The conditional-formatting rules on columns D, G, H, I, J were scoped
# to "..2:..49" (the old data extent). Re-scope each of them to include
# the freshly appended row 50.
function Extend-ConditionalFormatting($ws, $oldRange, $newRange) {
    $fcs = $ws.Range($oldRange).FormatConditions
    $count = $fcs.Count
    for ($i = 1; $i -le $count; $i++) {
        $fc = $fcs.Item($i)
        $fc.ModifyAppliesToRange($ws.Range($newRange))
    }
}

Extend-ConditionalFormatting $logs "D2:D49" "D2:D50"
Extend-ConditionalFormatting $logs "G2:G49" "G2:G50"
Extend-ConditionalFormatting $logs "H2:H49" "H2:H50"
Extend-ConditionalFormatting $logs "I2:I49" "I2:I50"
Extend-ConditionalFormatting $logs "J2:J49" "J2:J50"

# Dashboard tally for "Planning / Afspraak" goes from 13 to 14 to
# reflect the newly logged mail.
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Cells.Item(2, 2).Value = 14
